$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the slightly more precise timestamp for row 89 (A89)
$ws.Range("A89").Value = 44402.76951803009

# Add new row 90 of data
$ws.Range("A90").Value = 44403.76811139147
$ws.Range("B90").Value = 79931
$ws.Range("C90").Value = 67474
$ws.Range("D90").Value = 3579
$ws.Range("E90").Value = 2197
$ws.Range("F90").Value = 1598
$ws.Range("G90").Value = 20956
$ws.Range("H90").Value = 1571
$ws.Range("I90").Value = 896
$ws.Range("J90").Value = 199

# Apply same number format as other date cells in column A to the new cell
$ws.Range("A90").NumberFormat = $ws.Range("A89").NumberFormat
